$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price and Volume(1h) columns) per commit.
# Force the Price column to keep its text formatting so values like
# "258.85" or "2.011.68" are not reinterpreted as numbers/dates.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.156.68"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.011.68"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.85"
$ws.Range("E5").Value = "  +4.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  -1.83%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.21"
$ws.Range("E8").Value = "  -7.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  -3.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0769"
$ws.Range("E10").Value = "  -5.51%  "

$ws.Range("E11").Value = "  -2.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.39"
$ws.Range("E12").Value = "  -4.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.306.21"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.31"
$ws.Range("E14").Value = "  -3.01%  "

$ws.Range("E15").Value = "  -6.31%  "

$ws.Range("E16").Value = "  -4.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.009.33"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.045.94"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.34"
$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0835"
$ws.Range("E20").Value = "  -3.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "234.38"
$ws.Range("E21").Value = "  +1.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.09"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.83"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("E27").Value = "  -4.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.54"
$ws.Range("E28").Value = "  -1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.127"
$ws.Range("E29").Value = "  -7.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  -2.76%  "

$ws.Range("E31").Value = "  -2.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.58"
$ws.Range("E32").Value = "  -4.07%  "

$ws.Range("E33").Value = "  -4.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.41"
$ws.Range("E34").Value = "  -1.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("E35").Value = "  -6.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.54"
$ws.Range("E36").Value = "  -2.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  +0.74%  "

$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.53"
$ws.Range("E39").Value = "  +2.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.03"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.439.50"
$ws.Range("E42").Value = "  +4.27%  "

$ws.Range("E43").Value = "  -5.37%  "

$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.26"
$ws.Range("E45").Value = "  -2.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.54"
$ws.Range("E46").Value = "  -8.33%  "

$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.93"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.94"
$ws.Range("E49").Value = "  -6.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.198.05"
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("E51").Value = "  -8.06%  "
